$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1720.2858
$ws.Range("I6").Value = 1720.2858
$ws.Range("K6").Value = 5160.857400000001
$ws.Range("M6").Value = -5048.857400000001
$ws.Range("H9").Value = 146
$ws.Range("J9").Value = 79
$ws.Range("L9").Value = 79
$ws.Range("N9").Value = -417
$ws.Range("H38").Value = 37.375
$ws.Range("I38").Value = 37.375
$ws.Range("K38").Value = 112.125
$ws.Range("M38").Value = 259.875
$ws.Range("H40").Value = 5032.467
$ws.Range("I40").Value = 1993.5
$ws.Range("K40").Value = 1993.5
$ws.Range("M40").Value = -1818.5
$ws.Range("H43").Value = 5111.875
$ws.Range("I43").Value = 5724.5
$ws.Range("J43").Value = 4499.25
$ws.Range("K43").Value = 5724.5
$ws.Range("L43").Value = 4499.25
$ws.Range("M43").Value = -5655.5
$ws.Range("N43").Value = -4637.25
$ws.Range("H53").Value = 171.16667
$ws.Range("J53").Value = 120.75
$ws.Range("L53").Value = 120.75
$ws.Range("N53").Value = -1394.75
$ws.Range("H55").Value = 109.458336
$ws.Range("I55").Value = 115.416664
$ws.Range("J55").Value = 103.5
$ws.Range("K55").Value = 115.416664
$ws.Range("L55").Value = 103.5
$ws.Range("M55").Value = 98.583336
$ws.Range("N55").Value = -531.5
$ws.Range("H62").Value = 4243.6
$ws.Range("J62").Value = 4765.5
$ws.Range("L62").Value = 4765.5
$ws.Range("N62").Value = -6013.5
$ws.Range("H65").Value = 4243.6
$ws.Range("J65").Value = 4765.5
$ws.Range("L65").Value = 23827.5
$ws.Range("N65").Value = -30067.5
$ws.Range("H98").Value = 1702.3334
$ws.Range("I98").Value = 1592.7307
$ws.Range("K98").Value = 1592.7307
$ws.Range("M98").Value = -94.73070000000007
$ws.Range("H103").Value = 1899.875
$ws.Range("I103").Value = 1899.5
$ws.Range("J103").Value = 1900
$ws.Range("K103").Value = 5698.5
$ws.Range("L103").Value = 5700
$ws.Range("M103").Value = -5112.5
$ws.Range("N103").Value = -6872
$ws.Range("H106").Value = 5464.619
$ws.Range("I106").Value = 3309.9375
$ws.Range("K106").Value = 3309.9375
$ws.Range("M106").Value = -2678.9375
$ws.Range("H122").Value = 1702.3334
$ws.Range("I122").Value = 1592.7307
$ws.Range("K122").Value = 4778.1921
$ws.Range("M122").Value = -2328.1921
$ws.Range("H123").Value = 76890
$ws.Range("J123").Value = 76890
$ws.Range("L123").Value = 76890
$ws.Range("N123").Value = -86690
$ws.Range("H125").Value = 4029.2727
$ws.Range("J125").Value = 3983.75
$ws.Range("L125").Value = 35853.75
$ws.Range("N125").Value = -40773.75
$ws.Range("H138").Value = 3979.2856
$ws.Range("I138").Value = 3979.2856
$ws.Range("K138").Value = 11937.8568
$ws.Range("M138").Value = -6797.856800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1428.8462
$ws.Range("I2").Value = 1231.2222
$ws.Range("K2").Value = 1231.2222
$ws.Range("M2").Value = -1118.2222
$ws.Range("H32").Value = 913575.7
$ws.Range("J32").Value = 21839.8
$ws.Range("L32").Value = 21839.8
$ws.Range("N32").Value = -22413.8
$ws.Range("H45").Value = 1577.375
$ws.Range("I45").Value = 1163.8
$ws.Range("J45").Value = 2266.6667
$ws.Range("K45").Value = 1163.8
$ws.Range("L45").Value = 2266.6667
$ws.Range("M45").Value = -786.8
$ws.Range("N45").Value = -3020.6667
$ws.Range("H61").Value = 4547882
$ws.Range("I61").Value = 2534.2222
$ws.Range("J61").Value = 25001948
$ws.Range("K61").Value = 2534.2222
$ws.Range("L61").Value = 25001948
$ws.Range("M61").Value = -2322.2222
$ws.Range("N61").Value = -25002372
$ws.Range("H63").Value = 2411
$ws.Range("I63").Value = 2250
$ws.Range("K63").Value = 2250
$ws.Range("M63").Value = -1564
$ws.Range("H66").Value = 2411
$ws.Range("I66").Value = 2250
$ws.Range("K66").Value = 11250
$ws.Range("M66").Value = -7818
$ws.Range("H74").Value = 694722.7
$ws.Range("I74").Value = 739879.9
$ws.Range("K74").Value = 739879.9
$ws.Range("M74").Value = -739005.9
$ws.Range("H77").Value = 694722.7
$ws.Range("I77").Value = 739879.9
$ws.Range("K77").Value = 3699399.5
$ws.Range("M77").Value = -3695031.5
$ws.Range("H110").Value = 1931.6111
$ws.Range("I110").Value = 1857.1428
$ws.Range("J110").Value = 1979
$ws.Range("K110").Value = 1857.1428
$ws.Range("L110").Value = 1979
$ws.Range("M110").Value = 187.8571999999999
$ws.Range("N110").Value = -6069
$ws.Range("H116").Value = 1428.8462
$ws.Range("I116").Value = 1231.2222
$ws.Range("K116").Value = 1231.2222
$ws.Range("M116").Value = 1062.7778
$ws.Range("H122").Value = 1281.5
$ws.Range("I122").Value = 1237.3
$ws.Range("J122").Value = 1428.8334
$ws.Range("K122").Value = 3711.9
$ws.Range("L122").Value = 4286.5002
$ws.Range("M122").Value = -1261.9
$ws.Range("N122").Value = -9186.5002
$ws.Range("H132").Value = 5290.357
$ws.Range("I132").Value = 3103.1667
$ws.Range("K132").Value = 9309.500100000001
$ws.Range("M132").Value = -6779.500100000001
$ws.Range("H136").Value = 4547882
$ws.Range("I136").Value = 2534.2222
$ws.Range("J136").Value = 25001948
$ws.Range("K136").Value = 7602.6666
$ws.Range("L136").Value = 75005844
$ws.Range("M136").Value = -5052.6666
$ws.Range("N136").Value = -75010944

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1428.8462
$ws.Range("I3").Value = 1231.2222
$ws.Range("K3").Value = 1231.2222
$ws.Range("M3").Value = -1117.2222
$ws.Range("H20").Value = 24681.656
$ws.Range("J20").Value = 13677.143
$ws.Range("L20").Value = 13677.143
$ws.Range("N20").Value = -14171.143
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H86").Value = 3741.8823
$ws.Range("I86").Value = 3156.111
$ws.Range("K86").Value = 3156.111
$ws.Range("M86").Value = -2033.111
$ws.Range("H89").Value = 3741.8823
$ws.Range("I89").Value = 3156.111
$ws.Range("K89").Value = 15780.555
$ws.Range("M89").Value = -10164.555
$ws.Range("H105").Value = 6478.4287
$ws.Range("I105").Value = 1890
$ws.Range("K105").Value = 1890
$ws.Range("M105").Value = -143
$ws.Range("H107").Value = 1534.3793
$ws.Range("I107").Value = 1285
$ws.Range("J107").Value = 1767.1333
$ws.Range("K107").Value = 1285
$ws.Range("L107").Value = 1767.1333
$ws.Range("M107").Value = 635
$ws.Range("N107").Value = -5607.1333
$ws.Range("H134").Value = 2033623.1
$ws.Range("I134").Value = 1143.9231
$ws.Range("K134").Value = 3431.7693
$ws.Range("M134").Value = -896.7692999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 29500
$ws.Range("J51").Value = 29500
$ws.Range("L51").Value = 29500
$ws.Range("N51").Value = -30972
$ws.Range("H58").Value = 25060736
$ws.Range("I58").Value = 55562696
$ws.Range("K58").Value = 55562696
$ws.Range("M58").Value = -55562493
$ws.Range("H60").Value = 18951.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 18951.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 18951.5
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -19973.5
$ws.Range("H61").Value = 29500
$ws.Range("J61").Value = 29500
$ws.Range("L61").Value = 29500
$ws.Range("N61").Value = -30196
$ws.Range("H62").Value = 2926.5652
$ws.Range("I62").Value = 2795.762
$ws.Range("J62").Value = 4300
$ws.Range("K62").Value = 2795.762
$ws.Range("L62").Value = 4300
$ws.Range("M62").Value = -2171.762
$ws.Range("N62").Value = -5548
$ws.Range("H65").Value = 2926.5652
$ws.Range("I65").Value = 2795.762
$ws.Range("J65").Value = 4300
$ws.Range("K65").Value = 13978.81
$ws.Range("L65").Value = 21500
$ws.Range("M65").Value = -10858.81
$ws.Range("N65").Value = -27740
$ws.Range("H136").Value = 25060736
$ws.Range("I136").Value = 55562696
$ws.Range("K136").Value = 166688088
$ws.Range("M136").Value = -166685538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3276912.5
$ws.Range("I5").Value = 2747893.2
$ws.Range("K5").Value = 8243679.600000001
$ws.Range("M5").Value = -8243567.600000001
$ws.Range("H32").Value = 2004469.2
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2004469.2
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 6013407.6
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -6013973.6
$ws.Range("H64").Value = 7070.3105
$ws.Range("J64").Value = 7967.727
$ws.Range("L64").Value = 23903.181
$ws.Range("N64").Value = -24443.181
$ws.Range("H67").Value = 7070.3105
$ws.Range("J67").Value = 7967.727
$ws.Range("L67").Value = 23903.181
$ws.Range("N67").Value = -25775.181
$ws.Range("H87").Value = 14569.643
$ws.Range("I87").Value = 5425
$ws.Range("K87").Value = 16275
$ws.Range("M87").Value = -15027
$ws.Range("H90").Value = 14569.643
$ws.Range("I90").Value = 5425
$ws.Range("K90").Value = 48825
$ws.Range("M90").Value = -42585
$ws.Range("H131").Value = 4424.48
$ws.Range("I131").Value = 1166.9231
$ws.Range("J131").Value = 5569.027
$ws.Range("K131").Value = 3500.7693
$ws.Range("L131").Value = 16707.081
$ws.Range("M131").Value = 1539.2307
$ws.Range("N131").Value = -26787.081
$ws.Range("H135").Value = 3276912.5
$ws.Range("I135").Value = 2747893.2
$ws.Range("K135").Value = 24731038.8
$ws.Range("M135").Value = -24728503.8
$ws.Range("H136").Value = 11156.154
$ws.Range("I136").Value = 6206
$ws.Range("J136").Value = 14250
$ws.Range("K136").Value = 18618
$ws.Range("L136").Value = 42750
$ws.Range("M136").Value = -13518
$ws.Range("N136").Value = -52950

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8764.355
$ws.Range("I70").Value = 9292.962
$ws.Range("J70").Value = 8041
$ws.Range("K70").Value = 9292.962
$ws.Range("L70").Value = 8041
$ws.Range("M70").Value = -9022.962
$ws.Range("N70").Value = -8581
$ws.Range("H73").Value = 8764.355
$ws.Range("I73").Value = 9292.962
$ws.Range("J73").Value = 8041
$ws.Range("K73").Value = 9292.962
$ws.Range("L73").Value = 8041
$ws.Range("M73").Value = -8356.962
$ws.Range("N73").Value = -9913
$ws.Range("H80").Value = 1000
$ws.Range("J80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("N80").Value = -2996
$ws.Range("H83").Value = 1000
$ws.Range("J83").Value = 1000
$ws.Range("L83").Value = 5000
$ws.Range("N83").Value = -14984
$ws.Range("H107").Value = 3700
$ws.Range("I107").Value = 3266.6667
$ws.Range("K107").Value = 3266.6667
$ws.Range("M107").Value = -1346.6667
$ws.Range("H113").Value = 1293.4615
$ws.Range("I113").Value = 1319.8
$ws.Range("K113").Value = 1319.8
$ws.Range("M113").Value = 850.2
$ws.Range("H122").Value = 2447.3125
$ws.Range("I122").Value = 2500.3333
$ws.Range("J122").Value = 2161
$ws.Range("K122").Value = 7500.999899999999
$ws.Range("L122").Value = 6483
$ws.Range("M122").Value = -5050.999899999999
$ws.Range("N122").Value = -11383
$ws.Range("H140").Value = 92000
$ws.Range("J140").Value = 92000
$ws.Range("L140").Value = 92000
$ws.Range("N140").Value = -102360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3316.6667
$ws.Range("I7").Value = 2851.7144
$ws.Range("J7").Value = 4944
$ws.Range("K7").Value = 2851.7144
$ws.Range("L7").Value = 4944
$ws.Range("M7").Value = -2739.7144
$ws.Range("N7").Value = -5168
$ws.Range("H22").Value = 2535.0303
$ws.Range("I22").Value = 2101.125
$ws.Range("J22").Value = 2943.4119
$ws.Range("K22").Value = 2101.125
$ws.Range("L22").Value = 2943.4119
$ws.Range("M22").Value = -1806.125
$ws.Range("N22").Value = -3533.4119
$ws.Range("H27").Value = 2535.0303
$ws.Range("I27").Value = 2101.125
$ws.Range("J27").Value = 2943.4119
$ws.Range("K27").Value = 2101.125
$ws.Range("L27").Value = 2943.4119
$ws.Range("M27").Value = -1994.125
$ws.Range("N27").Value = -3157.4119
$ws.Range("H40").Value = 3365.9
$ws.Range("I40").Value = 2690.5625
$ws.Range("K40").Value = 2690.5625
$ws.Range("M40").Value = -2554.5625
$ws.Range("H42").Value = 59026.5
$ws.Range("I42").Value = 58025
$ws.Range("K42").Value = 58025
$ws.Range("M42").Value = -57462
$ws.Range("H46").Value = 4234.952
$ws.Range("J46").Value = 5398.6665
$ws.Range("L46").Value = 5398.6665
$ws.Range("N46").Value = -5774.6665
$ws.Range("H49").Value = 59026.5
$ws.Range("I49").Value = 58025
$ws.Range("K49").Value = 58025
$ws.Range("M49").Value = -57878
$ws.Range("H55").Value = 1618.6666
$ws.Range("I55").Value = 1445.5
$ws.Range("J55").Value = 1746.2632
$ws.Range("K55").Value = 1445.5
$ws.Range("L55").Value = 1746.2632
$ws.Range("M55").Value = -1272.5
$ws.Range("N55").Value = -2092.2632
$ws.Range("H68").Value = 2444.3333
$ws.Range("I68").Value = 2533.3333
$ws.Range("J68").Value = 1999.3334
$ws.Range("K68").Value = 2533.3333
$ws.Range("L68").Value = 1999.3334
$ws.Range("M68").Value = -1784.3333
$ws.Range("N68").Value = -3497.3334
$ws.Range("H71").Value = 2444.3333
$ws.Range("I71").Value = 2533.3333
$ws.Range("J71").Value = 1999.3334
$ws.Range("K71").Value = 12666.6665
$ws.Range("L71").Value = 9996.666999999999
$ws.Range("M71").Value = -8922.666499999999
$ws.Range("N71").Value = -17484.667
$ws.Range("H122").Value = 3038.423
$ws.Range("I122").Value = 2650.9
$ws.Range("K122").Value = 7952.700000000001
$ws.Range("M122").Value = -5502.700000000001
$ws.Range("H126").Value = 3316.6667
$ws.Range("I126").Value = 2851.7144
$ws.Range("J126").Value = 4944
$ws.Range("K126").Value = 8555.143199999999
$ws.Range("L126").Value = 14832
$ws.Range("M126").Value = -6085.143199999999
$ws.Range("N126").Value = -19772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18723.75
$ws.Range("I62").Value = 33263.332
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 33263.332
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -32639.332
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 18723.75
$ws.Range("I65").Value = 33263.332
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 166316.66
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -163196.66
$ws.Range("N65").Value = -56240
$ws.Range("H100").Value = 899.5
$ws.Range("I100").Value = 899.5
$ws.Range("K100").Value = 1799
$ws.Range("M100").Value = -1258
$ws.Range("H107").Value = 901.1875
$ws.Range("I107").Value = 801.46155
$ws.Range("J107").Value = 1333.3334
$ws.Range("K107").Value = 2404.38465
$ws.Range("L107").Value = 4000.0002
$ws.Range("M107").Value = -484.38465
$ws.Range("N107").Value = -7840.0002
$ws.Range("H113").Value = 1257.5116
$ws.Range("I113").Value = 995.36365
$ws.Range("J113").Value = 1532.1428
$ws.Range("K113").Value = 2986.09095
$ws.Range("L113").Value = 4596.428400000001
$ws.Range("M113").Value = -816.0909499999998
$ws.Range("N113").Value = -8936.428400000001
$ws.Range("H126").Value = 1850
$ws.Range("I126").Value = 1450
$ws.Range("K126").Value = 4350
$ws.Range("M126").Value = -1880
$ws.Range("H132").Value = 4388523.5
$ws.Range("I132").Value = 4507010.5
$ws.Range("K132").Value = 13521031.5
$ws.Range("M132").Value = -13518501.5
$ws.Range("H136").Value = 4450025
$ws.Range("I136").Value = 2387181.2
$ws.Range("J136").Value = 12701401
$ws.Range("K136").Value = 7161543.600000001
$ws.Range("L136").Value = 38104203
$ws.Range("M136").Value = -7158993.600000001
$ws.Range("N136").Value = -38109303
$ws.Range("H140").Value = 57990
$ws.Range("J140").Value = 57990
$ws.Range("L140").Value = 57990
$ws.Range("N140").Value = -68350
$ws.Range("H141").Value = 96690
$ws.Range("J141").Value = 96690
$ws.Range("L141").Value = 96690
$ws.Range("N141").Value = -107050
